$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (D) and 1h-volume-change (E) cells per the latest scrape.
# Cells whose new text still looks like a plain number (e.g. "326.93") are written
# with a leading apostrophe so Excel keeps them as text, matching the original
# inline-string cell type instead of silently converting them to numeric cells.
$ws.Range("D2").Value = "29.452.30"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "1.907.72"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'326.93"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4676"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("D9").Value = "'47.57"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "'0.08018"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("D11").Value = "'1.007"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("D12").Value = "'22.32"
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("D13").Value = "1.902.49"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "'5.934"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "'7.118"
$ws.Range("D16").Value = "'89.09"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("D18").Value = "'0.06591"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'0.00001028"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "'17.69"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "29.484.99"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").Value = "'5.526"
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("D24").Value = "'11.53"
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("D25").Value = "'2.209"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").Value = "2.155.58"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").Value = "'153.23"
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("D28").Value = "'19.77"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").Value = "'2.131"
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("D30").Value = "'5.707"
$ws.Range("E30").Value = "  +5.68%  "
$ws.Range("D31").Value = "'116.83"
$ws.Range("E31").Value = "  -1.66%  "
$ws.Range("D32").Value = "'1.072"
$ws.Range("E32").Value = "  +9.41%  "
$ws.Range("D33").Value = "'0.09492"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").Value = "'1.420"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D36").Value = "'5.379"
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("D37").Value = "'0.02256"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").Value = "'0.06078"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("D39").Value = "'8.358"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'1.171"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").Value = "'0.5867"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").Value = "'0.1835"
$ws.Range("D43").Value = "'10.11"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'2.452"
$ws.Range("E44").Value = "  +5.06%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.298"
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("D46").Value = "'0.07721"
$ws.Range("E46").Value = "  +9.93%  "
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("D48").Value = "'0.5537"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "'1.931"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("D50").Value = "'113.34"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").Value = "'0.2948"
$ws.Range("E51").Value = "  +6.21%  "
